# Cau 03 Va Cau 04 v2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Username) : written top-to-bottom, row 7 -> row 17 ---
$ws.Cells.Item(7, 1).Value = "user1"
$ws.Cells.Item(8, 1).Value = "user2"
$ws.Cells.Item(9, 1).Value = "user3"
$ws.Cells.Item(10, 1).Value = "user4"
$ws.Cells.Item(11, 1).Value = "user5"
$ws.Cells.Item(12, 1).Value = "user6"
$ws.Cells.Item(13, 1).Value = "user7"
$ws.Cells.Item(14, 1).Value = "user8"
$ws.Cells.Item(15, 1).Value = "user9"
$ws.Cells.Item(16, 1).Value = "user10"
$ws.Cells.Item(17, 1).Value = "user11"

# --- Column B (Password) : entered in this exact (non-linear) order ---
$ws.Cells.Item(7, 2).Value = "12321ádfasfdfsf"
$ws.Cells.Item(8, 2).Value = "adsfdsfsda"
$ws.Cells.Item(9, 2).Value = "ấdfsavdfv"
$ws.Cells.Item(10, 2).Value = "vxzvczx"
$ws.Cells.Item(12, 2).Value = "zcxvxczv"
$ws.Cells.Item(11, 2).Value = "agfdgfadg"
$ws.Cells.Item(14, 2).Value = "gadgfag"
$ws.Cells.Item(13, 2).Value = "fdgagfd"
$ws.Cells.Item(16, 2).Value = "àdgfda"
$ws.Cells.Item(15, 2).Value = "fadg"
$ws.Cells.Item(17, 2).Value = "sdfgdfs"

# --- Column C (Role) : numeric values, row 7 -> row 17 ---
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(17, 3).Value = 3

# Last row's password cell got left-aligned horizontally
$ws.Range("B17").HorizontalAlignment = -4131

# Final selection left on F11, matching the saved workbook state
$ws.Range("F11").Select()
